$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HOKCY")

$ws.Range("D8").Value = 4137200
$ws.Range("E8").Value = 3637900
$ws.Range("F8").Value = 3769600
$ws.Range("G8").Value = 4027400
$ws.Range("H8").Value = 3598200
$ws.Range("I8").Value = 3174900
$ws.Range("J8").Value = 2857000
$ws.Range("D9").Value = 1999000
$ws.Range("E9").Value = 1619100
$ws.Range("F9").Value = 1795900
$ws.Range("G9").Value = 2076300
$ws.Range("H9").Value = 1875400
$ws.Range("I9").Value = 1698400
$ws.Range("J9").Value = 1420700
$ws.Range("D10").Value = 2138200
$ws.Range("E10").Value = 2018800
$ws.Range("F10").Value = 1973700
$ws.Range("G10").Value = 1951100
$ws.Range("H10").Value = 1722800
$ws.Range("I10").Value = 1476500
$ws.Range("J10").Value = 1436300
$ws.Range("I12").Value = 15300
$ws.Range("D14").Value = -35300
$ws.Range("I14").Value = -84600
$ws.Range("D15").Value = 302500
$ws.Range("E15").Value = 283200
$ws.Range("F15").Value = 264400
$ws.Range("G15").Value = 248600
$ws.Range("H15").Value = 210100
$ws.Range("I15").Value = 186600
$ws.Range("J15").Value = 167000
$ws.Range("D17").Value = 3084800
$ws.Range("E17").Value = 2728300
$ws.Range("F17").Value = 2866300
$ws.Range("G17").Value = 3049900
$ws.Range("H17").Value = 2621900
$ws.Range("I17").Value = 2271100
$ws.Range("J17").Value = 2134100
$ws.Range("D18").Value = 1052400
$ws.Range("E18").Value = 909600
$ws.Range("F18").Value = 903300
$ws.Range("G18").Value = 977500
$ws.Range("H18").Value = 976400
$ws.Range("I18").Value = 903800
$ws.Range("J18").Value = 722900
$ws.Range("D20").Value = 521300
$ws.Range("E20").Value = 498400
$ws.Range("F20").Value = 502400
$ws.Range("G20").Value = 409500
$ws.Range("H20").Value = 340400
$ws.Range("I20").Value = 465600
$ws.Range("J20").Value = 400800
$ws.Range("D21").Value = 1876300
$ws.Range("E21").Value = 1691300
$ws.Range("F21").Value = 1670100
$ws.Range("G21").Value = 1635600
$ws.Range("H21").Value = 1526900
$ws.Range("I21").Value = 1556000
$ws.Range("J21").Value = "NA"
$ws.Range("D22").Value = 160100
$ws.Range("E22").Value = 153800
$ws.Range("F22").Value = 143800
$ws.Range("G22").Value = 129000
$ws.Range("H22").Value = 117900
$ws.Range("I22").Value = 110000
$ws.Range("J22").Value = 95800
$ws.Range("D23").Value = 1413600
$ws.Range("E23").Value = 1254200
$ws.Range("F23").Value = 1261900
$ws.Range("G23").Value = 1257900
$ws.Range("H23").Value = 1198800
$ws.Range("I23").Value = 1259300
$ws.Range("J23").Value = 1027900
$ws.Range("D24").Value = 222900
$ws.Range("E24").Value = 200800
$ws.Range("F24").Value = 220000
$ws.Range("G24").Value = 225700
$ws.Range("H24").Value = 210900
$ws.Range("I24").Value = 189100
$ws.Range("J24").Value = 171200
$ws.Range("D26").Value = 1190700
$ws.Range("E26").Value = 1053500
$ws.Range("F26").Value = 1042000
$ws.Range("G26").Value = 1032300
$ws.Range("H26").Value = 988000
$ws.Range("I26").Value = 1070200
$ws.Range("J26").Value = 856700
$ws.Range("D27").Value = 1047800
$ws.Range("E27").Value = 935100
$ws.Range("F27").Value = 930200
$ws.Range("G27").Value = 905600
$ws.Range("H27").Value = 873100
$ws.Range("I27").Value = 982400
$ws.Range("J27").Value = 783400
$ws.Range("D32").Value = -521300
$ws.Range("E32").Value = -498400
$ws.Range("F32").Value = -502400
$ws.Range("G32").Value = -409500
$ws.Range("H32").Value = -340400
$ws.Range("I32").Value = -465600
$ws.Range("J32").Value = -400800
$ws.Range("D33").Value = 1047800
$ws.Range("E33").Value = 935100
$ws.Range("F33").Value = 930200
$ws.Range("G33").Value = 905600
$ws.Range("H33").Value = 873100
$ws.Range("I33").Value = 982400
$ws.Range("J33").Value = 783400
$ws.Range("D35").Value = 1047800
$ws.Range("E35").Value = 935100
$ws.Range("F35").Value = 930200
$ws.Range("G35").Value = 905600
$ws.Range("H35").Value = 873100
$ws.Range("I35").Value = 982400
$ws.Range("J35").Value = 783400
$ws.Range("D41").Value = 1370500
$ws.Range("E41").Value = 1028800
$ws.Range("F41").Value = 1519200
$ws.Range("G41").Value = 1605800
$ws.Range("H41").Value = 1127300
$ws.Range("I41").Value = 2208300
$ws.Range("J41").Value = 1432100
$ws.Range("D42").Value = 284400
$ws.Range("E42").Value = 450400
$ws.Range("F42").Value = 170600
$ws.Range("G42").Value = 161600
$ws.Range("H42").Value = 248500
$ws.Range("I42").Value = 77500
$ws.Range("J42").Value = 102800
$ws.Range("D43").Value = 909200
$ws.Range("E43").Value = 788400
$ws.Range("F43").Value = 870600
$ws.Range("G43").Value = 826600
$ws.Range("H43").Value = 858900
$ws.Range("I43").Value = 677000
$ws.Range("J43").Value = 611600
$ws.Range("D44").Value = 328400
$ws.Range("E44").Value = 268800
$ws.Range("F44").Value = 291900
$ws.Range("G44").Value = 290900
$ws.Range("H44").Value = 303600
$ws.Range("I44").Value = 263300
$ws.Range("J44").Value = 206700
$ws.Range("D45").Value = 211300
$ws.Range("E45").Value = 160400
$ws.Range("F45").Value = 158300
$ws.Range("G45").Value = 254200
$ws.Range("H45").Value = 224700
$ws.Range("I45").Value = 190700
$ws.Range("J45").Value = 188800
$ws.Range("D46").Value = 3104000
$ws.Range("E46").Value = 2697000
$ws.Range("F46").Value = 3010600
$ws.Range("G46").Value = 3139100
$ws.Range("H46").Value = 2762900
$ws.Range("I46").Value = 2731000
$ws.Range("J46").Value = 2542100
$ws.Range("D47").Value = 4948100
$ws.Range("E47").Value = 4482200
$ws.Range("F47").Value = 4281400
$ws.Range("G47").Value = 3754500
$ws.Range("H47").Value = 3734200
$ws.Range("I47").Value = 3677800
$ws.Range("J47").Value = 3214600
$ws.Range("D48").Value = 7209200
$ws.Range("E48").Value = 6361700
$ws.Range("F48").Value = 6136200
$ws.Range("G48").Value = 6417700
$ws.Range("H48").Value = 5872300
$ws.Range("I48").Value = 3910000
$ws.Range("J48").Value = 4347100
$ws.Range("D49").Value = 749500
$ws.Range("E49").Value = 709900
$ws.Range("F49").Value = 741300
$ws.Range("G49").Value = 746300
$ws.Range("H49").Value = 669200
$ws.Range("I49").Value = 489900
$ws.Range("J49").Value = 437600
$ws.Range("D52").Value = 685200
$ws.Range("E52").Value = 621300
$ws.Range("F52").Value = 552100
$ws.Range("G52").Value = 517200
$ws.Range("H52").Value = 515100
$ws.Range("I52").Value = 470500
$ws.Range("J52").Value = 298100
$ws.Range("D54").Value = 16696000
$ws.Range("E54").Value = 14871900
$ws.Range("F54").Value = 14721600
$ws.Range("G54").Value = 14574800
$ws.Range("H54").Value = 13553700
$ws.Range("I54").Value = 12603600
$ws.Range("J54").Value = 10839500
$ws.Range("D57").Value = 1962800
$ws.Range("E57").Value = 1637400
$ws.Range("F57").Value = 1593500
$ws.Range("G57").Value = 1607700
$ws.Range("H57").Value = 1512000
$ws.Range("I57").Value = 1238500
$ws.Range("J57").Value = 1021900
$ws.Range("D58").Value = 2007300
$ws.Range("E58").Value = 758200
$ws.Range("F58").Value = 1237200
$ws.Range("G58").Value = 898100
$ws.Range("H58").Value = 792700
$ws.Range("I58").Value = 826900
$ws.Range("J58").Value = 537700
$ws.Range("D59").Value = 99800
$ws.Range("E59").Value = 94600
$ws.Range("F59").Value = 122200
$ws.Range("G59").Value = 129900
$ws.Range("H59").Value = 149100
$ws.Range("I59").Value = 132500
$ws.Range("J59").Value = 147800
$ws.Range("D60").Value = 4069900
$ws.Range("E60").Value = 2490200
$ws.Range("F60").Value = 2953000
$ws.Range("G60").Value = 2635600
$ws.Range("H60").Value = 2453800
$ws.Range("I60").Value = 2197800
$ws.Range("J60").Value = 1707500
$ws.Range("D61").Value = 2695800
$ws.Range("E61").Value = 3477300
$ws.Range("F61").Value = 2976300
$ws.Range("G61").Value = 3119100
$ws.Range("H61").Value = 3111300
$ws.Range("I61").Value = 3219100
$ws.Range("J61").Value = 2755200
$ws.Range("D62").Value = 981700
$ws.Range("E62").Value = 892000
$ws.Range("F62").Value = 879800
$ws.Range("G62").Value = 893400
$ws.Range("H62").Value = 807600
$ws.Range("I62").Value = 772600
$ws.Range("J62").Value = 474500
$ws.Range("D66").Value = 8696800
$ws.Range("E66").Value = 7701700
$ws.Range("F66").Value = 7690300
$ws.Range("G66").Value = 7570600
$ws.Range("H66").Value = 7201100
$ws.Range("I66").Value = 6845400
$ws.Range("J66").Value = 5542100
$ws.Range("D72").Value = 6772900
$ws.Range("E72").Value = 6299200
$ws.Range("F72").Value = 5906400
$ws.Range("G72").Value = 5518300
$ws.Range("H72").Value = 5097900
$ws.Range("I72").Value = 4600100
$ws.Range("J72").Value = 4154200
$ws.Range("D76").Value = 7999200
$ws.Range("E76").Value = 7170300
$ws.Range("F76").Value = 7031300
$ws.Range("G76").Value = 7004200
$ws.Range("H76").Value = 6352600
$ws.Range("I76").Value = 5758300
$ws.Range("J76").Value = 5297400
$ws.Range("D81").Value = 1047800
$ws.Range("E81").Value = 935100
$ws.Range("F81").Value = 930200
$ws.Range("G81").Value = 905600
$ws.Range("H81").Value = 873100
$ws.Range("I81").Value = 982400
$ws.Range("J81").Value = 783400
$ws.Range("D83").Value = 302500
$ws.Range("E83").Value = 283200
$ws.Range("F83").Value = 264400
$ws.Range("G83").Value = 248600
$ws.Range("H83").Value = 210100
$ws.Range("I83").Value = 186600
$ws.Range("J83").Value = "NA"
$ws.Range("D89").Value = 1085900
$ws.Range("E89").Value = 1078900
$ws.Range("F89").Value = 1054300
$ws.Range("G89").Value = 1042000
$ws.Range("H89").Value = 890700
$ws.Range("I89").Value = 849100
$ws.Range("J89").Value = 631500
$ws.Range("D91").Value = -764900
$ws.Range("E91").Value = -771600
$ws.Range("F91").Value = -768100
$ws.Range("G91").Value = -796200
$ws.Range("H91").Value = -647900
$ws.Range("I91").Value = -740900
$ws.Range("J91").Value = -583300
$ws.Range("D94").Value = -212800
$ws.Range("E94").Value = -824000
$ws.Range("F94").Value = -603700
$ws.Range("G94").Value = -360900
$ws.Range("H94").Value = -730400
$ws.Range("I94").Value = -802900
$ws.Range("J94").Value = "NA"
$ws.Range("D96").Value = -586400
$ws.Range("E96").Value = -533100
$ws.Range("F96").Value = -484800
$ws.Range("G96").Value = -440800
$ws.Range("H96").Value = -400800
$ws.Range("I96").Value = -540500
$ws.Range("J96").Value = -331200
$ws.Range("D100").Value = -552800
$ws.Range("E100").Value = -722700
$ws.Range("F100").Value = -506600
$ws.Range("G100").Value = -184000
$ws.Range("H100").Value = -595600
$ws.Range("I100").Value = 73800
$ws.Range("J100").Value = "NA"
$ws.Range("D101").Value = 21400
$ws.Range("E101").Value = -22600
$ws.Range("F101").Value = -30600
$ws.Range("G101").Value = -18600
$ws.Range("H101").Value = 10100
$ws.Range("I101").Value = 200
$ws.Range("J101").Value = "NA"
$ws.Range("D102").Value = 341700
$ws.Range("E102").Value = -490400
$ws.Range("F102").Value = -86600
$ws.Range("G102").Value = 478500
$ws.Range("H102").Value = -425200
$ws.Range("I102").Value = 120300
$ws.Range("J102").Value = 196900
